# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# The Tight End (TE) roster gained a new player, J.Sprinkle, whose
# Week 16 simulated stats are all zero to start.
$wsTE = $wb.Worksheets.Item("TE")

$wsTE.Range("A5").Value = "J.Sprinkle"
$wsTE.Range("B5:J5").Value = 0

# TE becomes the active sheet (was RB before), with J6 selected.
$wsTE.Activate()
$wsTE.Range("J6").Select()
